$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New MLP (neural network) results appended as rows 16 and 17, matching
# the unstyled formatting already used by the preceding appended rows
# (11-15), not the styled header block (rows 1-10).

$ws.Cells.Item(16, 1).Value = "{'activation': 'relu', 'alpha': 1e-06, 'beta_1': 0.9, 'hidden_layer_sizes': 10, 'learning_rate': 'constant', 'learning_rate_init': 0.01, 'max_iter': 3000, 'momentum': 0.9, 'power_t': 0.5, 'random_state': 6, 'solver': 'lbfgs'}"
$ws.Cells.Item(16, 2).Value = 0.992
$ws.Cells.Item(16, 3).Value = 0.783

$ws.Cells.Item(17, 1).Value = "{'activation': 'relu', 'alpha': 1e-06, 'beta_1': 0.9, 'hidden_layer_sizes': 10, 'learning_rate': 'constant', 'learning_rate_init': 0.1, 'max_iter': 500, 'momentum': 0.9, 'power_t': 0.5, 'random_state': 6, 'solver': 'lbfgs'}"
$ws.Cells.Item(17, 2).Value = 0.992
$ws.Cells.Item(17, 3).Value = 0.783

# Row 15 carries no explicit cell style (style index 0); copy that
# formatting onto the two new rows so they don't inherit the column's
# default styled format (index 3) the way a plain .Value write would.
$ws.Range("A15:C15").Copy() | Out-Null
$ws.Range("A16:C17").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
